# Dividend Calculation workbook update
# - Suzie's Roth IRA dividend for October ("Yearly" sheet, F12) increases from 0 to 7.88,
#   which cascades through the dependent SUM formulas (G12, F15, G15).
# - The "All Time" sheet keeps a manually-copied snapshot of the Yearly totals for each
#   year (F/G/H columns), so the 2016 401K total (H7) must be updated to match the new
#   Yearly F15 total; its dependent SUM formulas (I7, H46, I46) cascade automatically.
# - The active sheet/selection also moved from "Yearly" to "All Time".

$wb = $excel.ActiveWorkbook

$wsYearly  = $wb.Worksheets.Item("Yearly")
$wsAllTime = $wb.Worksheets.Item("All Time")

# Suzie's Roth IRA dividend for October (Yearly!F12): 0 -> 7.88
$wsYearly.Range("F12").Value = 7.88

# All Time!H7 mirrors Yearly!F15 (the 401K/"Suzie's Roth IRA" yearly total for 2016)
$wsAllTime.Range("H7").Value = 137.80000000000001

# Update the view/selection state:
#  - "Yearly" is no longer the selected tab; its selection moves to F15.
$wsYearly.Range("F15").Select()

#  - "All Time" becomes the active/selected tab, scrolled to row 13, selecting I13.
$wsAllTime.Activate()
$wsAllTime.Range("I13").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
